$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.546166777610779
$ws.Range("B1").Value = 2.447568416595459
$ws.Range("C1").Value = 4.516322135925293
$ws.Range("D1").Value = 4.358202457427979
$ws.Range("E1").Value = 1.409067153930664
